$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename (row 1) ---
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# --- Title-case fix for Spanish particles (de/del/la/las/el/los/y) in state/municipality names ---
$ws.Range('B5').Value = 'Pabellón De Arteaga'
$ws.Range('B6').Value = 'Rincón De Romos'
$ws.Range('B24').Value = 'Benemérito De Las Américas'
$ws.Range('B49').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B55').Value = 'Salto De Agua'
$ws.Range('B75').Value = 'Coyame Del Sotol'
$ws.Range('B80').Value = 'Hidalgo Del Parral'
$ws.Range('B86').Value = 'San Francisco Del Oro'
$ws.Range('B89').Value = 'Valle De Zaragoza'
$ws.Range('B103').Value = 'San Juan De Sabinas'
$ws.Range('B110').Value = 'Villa De Álvarez'
$ws.Range('A112').Value = 'Ciudad De México'
$ws.Range('B116').Value = 'Cuajimalpa De Morelos'
$ws.Range('B132').Value = 'Nombre De Dios'
$ws.Range('B133').Value = 'Pánuco De Coronado'
$ws.Range('B136').Value = 'San Juan Del Río'
$ws.Range('B137').Value = 'San Pedro Del Gallo'
$ws.Range('A141').Value = 'Estado De México'
$ws.Range('B141').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B150').Value = 'Coacalco De Berriozábal'
$ws.Range('B153').Value = 'Ecatepec De Morelos'
$ws.Range('B158').Value = 'Naucalpan De Juárez'
$ws.Range('B161').Value = 'San Felipe Del Progreso'
$ws.Range('B166').Value = 'Tenango Del Aire'
$ws.Range('B167').Value = 'Tenango Del Valle'
$ws.Range('B169').Value = 'Tlalnepantla De Baz'
$ws.Range('B172').Value = 'Valle De Bravo'
$ws.Range('B173').Value = 'Valle De Chalco Solidaridad'
$ws.Range('B174').Value = 'Villa Del Carbón'
$ws.Range('B180').Value = 'San Miguel De Allende'
$ws.Range('B181').Value = 'Apaseo El Alto'
$ws.Range('B182').Value = 'Apaseo El Grande'
$ws.Range('B189').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B192').Value = 'Jaral Del Progreso'
$ws.Range('B199').Value = 'Purísima Del Rincón'
$ws.Range('B202').Value = 'San Diego De La Unión'
$ws.Range('B205').Value = 'San Luis De La Paz'
$ws.Range('B206').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B207').Value = 'Silao De La Victoria'
$ws.Range('B211').Value = 'Valle De Santiago'
$ws.Range('B217').Value = 'Acapulco De Juárez'
$ws.Range('B219').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B220').Value = 'Alcozauca De Guerrero'
$ws.Range('B222').Value = 'Atenango Del Río'
$ws.Range('B223').Value = 'Atoyac De Álvarez'
$ws.Range('B224').Value = 'Ayutla De Los Libres'
$ws.Range('B227').Value = 'Buenavista De Cuéllar'
$ws.Range('B228').Value = 'Chilapa De Álvarez'
$ws.Range('B229').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B232').Value = 'Coyuca De Benítez'
$ws.Range('B233').Value = 'Coyuca De Catalán'
$ws.Range('B236').Value = 'Cuetzala Del Progreso'
$ws.Range('B237').Value = 'Cutzamala De Pinzón'
$ws.Range('B242').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B243').Value = 'Iguala De La Independencia'
$ws.Range('B244').Value = 'Zihuatanejo De Azueta'
$ws.Range('B255').Value = 'Taxco De Alarcón'
$ws.Range('B257').Value = 'Técpan De Galeana'
$ws.Range('B259').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B260').Value = 'Tixtla De Guerrero'
$ws.Range('B263').Value = 'Tlapa De Comonfort'
$ws.Range('B275').Value = 'Atotonilco El Grande'
$ws.Range('B279').Value = 'Cuautepec De Hinojosa'
$ws.Range('B282').Value = 'Huejutla De Reyes'
$ws.Range('B285').Value = 'Jacala De Ledezma'
$ws.Range('B292').Value = 'Mixquiahuala De Juárez'
$ws.Range('B294').Value = 'Nopala De Villagrán'
$ws.Range('B295').Value = 'Pachuca De Soto'
$ws.Range('B297').Value = 'Progreso De Obregón'
$ws.Range('B300').Value = 'Santiago De Anaya'
$ws.Range('B303').Value = 'Tenango De Doria'
$ws.Range('B305').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B307').Value = 'Tezontepec De Aldama'
$ws.Range('B312').Value = 'Tula De Allende'
$ws.Range('B313').Value = 'Tulancingo De Bravo'
$ws.Range('B314').Value = 'Zacualtipán De Ángeles'
$ws.Range('B319').Value = 'Atotonilco El Alto'
$ws.Range('B320').Value = 'Autlán De Navarro'
$ws.Range('B326').Value = 'Encarnación De Díaz'
$ws.Range('B329').Value = 'Jilotlán De Los Dolores'
$ws.Range('B331').Value = 'Lagos De Moreno'
$ws.Range('B333').Value = 'Ojuelos De Jalisco'
$ws.Range('B336').Value = 'Santa María De Los Ángeles'
$ws.Range('B339').Value = 'Tepatitlán De Morelos'
$ws.Range('B341').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B345').Value = 'Yahualica De González Gallo'
$ws.Range('B381').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B398').Value = 'Coatlán Del Río'
$ws.Range('B405').Value = 'Tlaltizapán De Zapata'
$ws.Range('B409').Value = 'Santa María Del Oro'
$ws.Range('B421').Value = 'Mier Y Noriega'
$ws.Range('B426').Value = 'San Nicolás De Los Garza'
$ws.Range('B430').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B434').Value = 'Cuilápam De Guerrero'
$ws.Range('B435').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B436').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B437').Value = 'Ixtlán De Juárez'
$ws.Range('B438').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B440').Value = 'Mariscala De Juárez'
$ws.Range('B442').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B443').Value = 'Oaxaca De Juárez'
$ws.Range('B444').Value = 'Ocotlán De Morelos'
$ws.Range('B445').Value = 'Putla Villa De Guerrero'
$ws.Range('B486').Value = 'Santo Domingo De Morelos'
$ws.Range('B489').Value = 'Tataltepec De Valdés'
$ws.Range('B490').Value = 'Teotitlán De Flores Magón'
$ws.Range('B491').Value = 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'
$ws.Range('B492').Value = 'Tlacolula De Matamoros'
$ws.Range('B495').Value = 'Zimatlán De Álvarez'
$ws.Range('B508').Value = 'Huehuetlán El Chico'
$ws.Range('B511').Value = 'Huitzilan De Serdán'
$ws.Range('B513').Value = 'Izúcar De Matamoros'
$ws.Range('B518').Value = 'Los Reyes De Juárez'
$ws.Range('B521').Value = 'Palmar De Bravo'
$ws.Range('B529').Value = 'Tepexi De Rodríguez'
$ws.Range('B530').Value = 'Tepeyahualco De Cuauhtémoc'
$ws.Range('B531').Value = 'Tetela De Ocampo'
$ws.Range('B533').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B538').Value = 'Tuzamapan De Galeana'
$ws.Range('B547').Value = 'Amealco De Bonfil'
$ws.Range('B549').Value = 'Cadereyta De Montes'
$ws.Range('B553').Value = 'Jalpan De Serra'
$ws.Range('B554').Value = 'Landa De Matamoros'
$ws.Range('B556').Value = 'Pinal De Amoles'
$ws.Range('B558').Value = 'San Juan Del Río'
$ws.Range('B568').Value = 'Axtla De Terrazas'
$ws.Range('B573').Value = 'Ciudad Del Maíz'
$ws.Range('B583').Value = 'Mexquitic De Carmona'
$ws.Range('B592').Value = 'Santa María Del Río'
$ws.Range('B594').Value = 'Soledad De Graciano Sánchez'
$ws.Range('B599').Value = 'Tanquián De Escobedo'
$ws.Range('B601').Value = 'Villa De Arista'
$ws.Range('B602').Value = 'Villa De Guadalupe'
$ws.Range('B603').Value = 'Villa De Ramos'
$ws.Range('B604').Value = 'Villa De Reyes'
$ws.Range('B655').Value = 'Soto La Marina'
$ws.Range('B663').Value = 'Nanacamilpa De Mariano Arista'
$ws.Range('B667').Value = 'Tetla De La Solidaridad'
$ws.Range('B678').Value = 'Amatlán De Los Reyes'
$ws.Range('B684').Value = 'Camarón De Tejeda'
$ws.Range('B692').Value = 'Cosamaloapan De Carpio'
$ws.Range('B701').Value = 'Ixhuatlán De Madero'
$ws.Range('B702').Value = 'Ixhuatlán Del Café'
$ws.Range('B703').Value = 'Ixhuatlán Del Sureste'
$ws.Range('B711').Value = 'Lerdo De Tejada'
$ws.Range('B712').Value = 'Martínez De La Torre'
$ws.Range('B719').Value = 'Ozuluama De Mascareñas'
$ws.Range('B722').Value = 'Paso Del Macho'
$ws.Range('B726').Value = 'Poza Rica De Hidalgo'
$ws.Range('B732').Value = 'Sayula De Alemán'
$ws.Range('B733').Value = 'Soledad De Doblado'
$ws.Range('B752').Value = 'Vega De Alatorre'
$ws.Range('B779').Value = 'Nochistlán De Mejía'
$ws.Range('B786').Value = 'Tlaltenango De Sánchez Román'

# --- Floating point precision corrections on pct_matriculas column ---
$ws.Range('D6').Value = 0.0009074410163339384
$ws.Range('D16').Value = 0.0009074410163339384
$ws.Range('D26').Value = 0.0009074410163339384
$ws.Range('D31').Value = 0.0009074410163339384
$ws.Range('D35').Value = 0.0009074410163339384
$ws.Range('D40').Value = 0.0009074410163339384
$ws.Range('D43').Value = 0.0009074410163339384
$ws.Range('D50').Value = 0.0009074410163339384
$ws.Range('D61').Value = 0.0009074410163339384
$ws.Range('D76').Value = 0.0009074410163339384
$ws.Range('D79').Value = 0.0009074410163339384
$ws.Range('D82').Value = 0.0009074410163339384
$ws.Range('D102').Value = 0.0009074410163339384
$ws.Range('D103').Value = 0.0009074410163339384
$ws.Range('D111').Value = 0.0009074410163339384
$ws.Range('D129').Value = 0.0009074410163339384
$ws.Range('D131').Value = 0.0009074410163339384
$ws.Range('D133').Value = 0.0009074410163339384
$ws.Range('D138').Value = 0.0009074410163339384
$ws.Range('D152').Value = 0.0009074410163339384
$ws.Range('D187').Value = 0.0009074410163339384
$ws.Range('D190').Value = 0.0009074410163339384
$ws.Range('D207').Value = 0.0009074410163339384
$ws.Range('D208').Value = 0.0009074410163339384
$ws.Range('D218').Value = 0.0009074410163339384
$ws.Range('D234').Value = 0.0009074410163339384
$ws.Range('D238').Value = 0.0009074410163339384
$ws.Range('D242').Value = 0.0009074410163339384
$ws.Range('D245').Value = 0.0009074410163339384
$ws.Range('D248').Value = 0.0009074410163339384
$ws.Range('D253').Value = 0.0009074410163339384
$ws.Range('D261').Value = 0.0009074410163339384
$ws.Range('D266').Value = 0.0009074410163339384
$ws.Range('D279').Value = 0.0009074410163339384
$ws.Range('D310').Value = 0.0009074410163339384
$ws.Range('D317').Value = 0.0009074410163339384
$ws.Range('D345').Value = 0.0009074410163339384
$ws.Range('D348').Value = 0.0009074410163339384
$ws.Range('D351').Value = 0.0009074410163339384
$ws.Range('D358').Value = 0.0009074410163339384
$ws.Range('D364').Value = 0.0009074410163339384
$ws.Range('D367').Value = 0.0009074410163339384
$ws.Range('D388').Value = 0.0009074410163339384
$ws.Range('D391').Value = 0.0009074410163339384
$ws.Range('D399').Value = 0.0009074410163339384
$ws.Range('D405').Value = 0.0009074410163339384
$ws.Range('D426').Value = 0.0009074410163339384
$ws.Range('D430').Value = 0.0009074410163339384
$ws.Range('D431').Value = 0.0009074410163339384
$ws.Range('D445').Value = 0.0009074410163339384
$ws.Range('D457').Value = 0.0009074410163339384
$ws.Range('D479').Value = 0.0009074410163339384
$ws.Range('D491').Value = 0.0009074410163339384
$ws.Range('D501').Value = 0.0009074410163339384
$ws.Range('D507').Value = 0.0009074410163339384
$ws.Range('D508').Value = 0.0009074410163339384
$ws.Range('D509').Value = 0.009679370840895344
$ws.Range('D512').Value = 0.0009074410163339384
$ws.Range('D515').Value = 0.0009074410163339384
$ws.Range('D526').Value = 0.0009074410163339384
$ws.Range('D539').Value = 0.00998185117967332
$ws.Range('D580').Value = 0.0009074410163339384
$ws.Range('D581').Value = 0.0009074410163339384
$ws.Range('D587').Value = 0.0009074410163339384
$ws.Range('D588').Value = 0.009679370840895344
$ws.Range('D595').Value = 0.0009074410163339384
$ws.Range('D598').Value = 0.0009074410163339384
$ws.Range('D599').Value = 0.0009074410163339384
$ws.Range('D604').Value = 0.0009074410163339384
$ws.Range('D621').Value = 0.0009074410163339384
$ws.Range('D627').Value = 0.0009074410163339384
$ws.Range('D638').Value = 0.0009074410163339384
$ws.Range('D652').Value = 0.0009074410163339384
$ws.Range('D656').Value = 0.009679370840895344
$ws.Range('D665').Value = 0.0009074410163339384
$ws.Range('D668').Value = 0.0009074410163339384
$ws.Range('D679').Value = 0.0009074410163339384
$ws.Range('D688').Value = 0.0009074410163339384
$ws.Range('D699').Value = 0.0009074410163339384
$ws.Range('D706').Value = 0.0009074410163339384
$ws.Range('D708').Value = 0.0009074410163339384
$ws.Range('D715').Value = 0.0009074410163339384
$ws.Range('D743').Value = 0.0009074410163339384
$ws.Range('D745').Value = 0.0009074410163339384
$ws.Range('D748').Value = 0.0009074410163339384
$ws.Range('D749').Value = 0.0009074410163339384
$ws.Range('D770').Value = 0.0009074410163339384
$ws.Range('D786').Value = 0.0009074410163339384

# --- Remove trailing footnote rows (795-799), shrinking used range to A1:D793 ---
$ws.Range("A795:D799").EntireRow.Delete()
